# Update "Dados BIBI" retention metrics per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: num_customers 26 -> 27 (cohort_size stays 2654); retention_rate recalculated
$ws.Range("C22").Value = 27
$ws.Range("E22").Value = 27 / 2654

# Row 34: num_customers 74 -> 76 (cohort_size stays 2256); retention_rate recalculated
$ws.Range("C34").Value = 76
$ws.Range("E34").Value = 76 / 2256

# Row 37: num_customers and cohort_size both 750 -> 755 (retention_rate stays 1)
$ws.Range("C37").Value = 755
$ws.Range("D37").Value = 755
